$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - TROXEL 1001: fill in Layout/Asthetic/Chairs scores and mark building-summary row complete
$ws.Range("D3").Value = 7.75
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 5.5
$ws.Range("M3").Value = 1

# Row 5 - KILDEE 0125
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("M5").Value = 1

# Row 22 - SICTR 2221 (scores already present) - mark building-summary row complete
$ws.Range("M22").Value = 1

# Row 25 - PHYSICS 0003 - mark building-summary row complete
$ws.Range("M25").Value = 1

# Row 28 - HORT 0118
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = 4

# Row 29 - ROSS H 0124 - mark building-summary row complete
$ws.Range("M29").Value = 1

# Row 31 - BESSEY 0210
$ws.Range("D31").Value = 7.5
$ws.Range("E31").Value = 6
$ws.Range("F31").Value = 2

# Row 38 - HAMILTN 0169
$ws.Range("D38").Value = 2.5
$ws.Range("E38").Value = 3
$ws.Range("F38").Value = 4

# Update the active selection / scroll position to reflect the latest edit location
$ws.Range("Q16").Select()
